# Update the Turkish fantasy-basketball roster ("LOS Galacticos") sheet:
#   - swap out two players (Dennis Schröder, Bobby Portis) for Draymond Green
#   - re-sort the roster into its new order
#   - the table shrinks from 18 players to 17, so the old last row is removed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old roster body first so stale shared strings (the two
# outgoing players) are dropped and the new table is written fresh.
$ws.Range("A2:C19").ClearContents()

$players = @(
    @("Anfernee Simons",        "PG,SG",    "Portland Trail Blazers"),
    @("Paul George",            "SG,SF,PF", "Philadelphia 76ers"),
    @("Fred VanVleet",          "PG",       "Houston Rockets"),
    @("James Harden",           "PG,SG",    "LA Clippers"),
    @("Dillon Brooks",          "SG,SF",    "Houston Rockets"),
    @("Jaren Jackson Jr.",      "PF,C",     "Memphis Grizzlies"),
    @("Anthony Edwards",        "SG,SF",    "Minnesota Timberwolves"),
    @("Amen Thompson",          "SG,SF",    "Houston Rockets"),
    @("Ivica Zubac",            "C",        "LA Clippers"),
    @("Draymond Green",         "PF,C",     "Golden State Warriors"),
    @("Bilal Coulibaly",        "SG,SF",    "Washington Wizards"),
    @("Andrew Wiggins",         "SF,PF",    "Golden State Warriors"),
    @("Zion Williamson",        "PF,C",     "New Orleans Pelicans"),
    @("Jayson Tatum",           "SF,PF",    "Boston Celtics"),
    @("Giannis Antetokounmpo",  "PF,C",     "Milwaukee Bucks"),
    @("Jonathan Kuminga",       "SF,PF",    "Golden State Warriors"),
    @("Jerami Grant",           "SF,PF",    "Portland Trail Blazers")
)

$row = 2
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p[0]
    $ws.Cells.Item($row, 2).Value = $p[1]
    $ws.Cells.Item($row, 3).Value = $p[2]
    $row++
}

# The roster used to have 18 players (through row 19); now it only has 17
# (through row 18), so drop the now-unused trailing row.
$ws.Rows.Item(19).Delete()
